$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.6707215716073449
$ws.Range("C2").Value = 0.1579930342821001
$ws.Range("D2").Value = 0.01294316665145345
$ws.Range("E2").Value = 0.1264417977370229
$ws.Range("F2").Value = 0.5538072680917523
$ws.Range("M2").Value = 0.3099252839110775
$ws.Range("N2").Value = 0.9722120263364573
$ws.Range("O2").Value = 1.797013954961699
$ws.Range("B3").Value = 0.5892113537391026
$ws.Range("C3").Value = 0.1408492153979353
$ws.Range("D3").Value = 0.01180183988952876
$ws.Range("E3").Value = 0.1199290195624414
$ws.Range("F3").Value = 0.5440222806960406
$ws.Range("M3").Value = 0.2762130393736228
$ws.Range("N3").Value = 0.9838772478455269
$ws.Range("O3").Value = 1.77823950383555
$ws.Range("B4").Value = 0.5391197371740475
$ws.Range("C4").Value = 0.1302612413019233
$ws.Range("D4").Value = 0.01109670045017097
$ws.Range("E4").Value = 0.1160393710671315
$ws.Range("F4").Value = 0.5384551950288596
$ws.Range("M4").Value = 0.2555826412098341
$ws.Range("N4").Value = 0.9914898772992551
$ws.Range("O4").Value = 1.768167780096036
$ws.Range("B5").Value = 0.5186967731642085
$ws.Range("C5").Value = 0.1259312547789762
$ws.Range("D5").Value = 0.01080827466421042
$ws.Range("E5").Value = 0.11448151404214
$ws.Range("F5").Value = 0.5362971772272331
$ws.Range("M5").Value = 0.247192905060615
$ws.Range("N5").Value = 0.994705237091658
$ws.Range("O5").Value = 1.764428586231787
$ws.Range("B6").Value = 0.5153049630227997
$ws.Range("C6").Value = 0.1252113452448214
$ws.Range("D6").Value = 0.01076031742006123
$ws.Range("E6").Value = 0.1142244687086773
$ws.Range("F6").Value = 0.5359455132181807
$ws.Range("M6").Value = 0.2458008404967487
$ws.Range("N6").Value = 0.9952459781507166
$ws.Range("O6").Value = 1.763829719546123
$ws.Range("B7").Value = 0.538844346213267
$ws.Range("C7").Value = 0.1302029072919879
$ws.Range("D7").Value = 0.01109281496717074
$ws.Range("E7").Value = 0.1160182514466896
$ws.Range("F7").Value = 0.5384256437198047
$ws.Range("M7").Value = 0.2554694242300144
$ws.Range("N7").Value = 0.9915327826713458
$ws.Range("O7").Value = 1.768115874978065
$ws.Range("B8").Value = 0.6426264640834063
$ws.Range("C8").Value = 0.1520947386414093
$ws.Range("D8").Value = 0.01255055418344142
$ws.Range("E8").Value = 0.124173364358235
$ws.Range("F8").Value = 0.5503417257166134
$ws.Range("M8").Value = 0.298286930748425
$ws.Range("N8").Value = 0.976140745274833
$ws.Range("O8").Value = 1.790237761937334
$ws.Range("B9").Value = 0.8457665818863234
$ws.Range("C9").Value = 0.1945294744807882
$ws.Range("D9").Value = 0.01537380291151891
$ws.Range("E9").Value = 0.1410447569340363
$ws.Range("F9").Value = 0.5772216663375218
$ws.Range("M9").Value = 0.3828063636020644
$ws.Range("N9").Value = 0.9495302693901948
$ws.Range("O9").Value = 1.845219442191166
$ws.Range("B10").Value = 0.9947626019351787
$ws.Range("C10").Value = 0.2253989016726052
$ws.Range("D10").Value = 0.01742562688302485
$ws.Range("E10").Value = 0.1539949508719047
$ws.Range("F10").Value = 0.5991339525288311
$ws.Range("M10").Value = 0.445256160204039
$ws.Range("N10").Value = 0.932159991699514
$ws.Range("O10").Value = 1.892761987288566
$ws.Range("B11").Value = 1.062486596629753
$ws.Range("C11").Value = 0.2393744685245451
$ws.Range("D11").Value = 0.01835401909765721
$ws.Range("E11").Value = 0.1600109248144221
$ws.Range("F11").Value = 0.6095771854180612
$ws.Range("M11").Value = 0.4737466196331752
$ws.Range("N11").Value = 0.9247317984629007
$ws.Range("O11").Value = 1.91595900815463
$ws.Range("B12").Value = 1.088123353767685
$ws.Range("C12").Value = 0.2446568546795334
$ws.Range("D12").Value = 0.01870484080907175
$ws.Range("E12").Value = 0.162307279819899
$ws.Range("F12").Value = 0.6136004393797805
$ws.Range("M12").Value = 0.4845471258391427
$ws.Range("N12").Value = 0.9219871133193749
$ws.Range("O12").Value = 1.924970030011423
$ws.Range("B13").Value = 1.082602427868153
$ws.Range("C13").Value = 0.2435196400683708
$ws.Range("D13").Value = 0.01862931835459136
$ws.Range("E13").Value = 0.1618119031397711
$ws.Range("F13").Value = 0.6127309027531567
$ws.Range("M13").Value = 0.4822205187690258
$ws.Range("N13").Value = 0.9225751950132661
$ws.Range("O13").Value = 1.923019240194549
$ws.Range("B14").Value = 1.064595930776704
$ws.Range("C14").Value = 0.2398092523775119
$ws.Range("D14").Value = 0.01838289637464641
$ws.Range("E14").Value = 0.1601994802296076
$ws.Range("F14").Value = 0.609906803860909
$ws.Range("M14").Value = 0.474634947338231
$ws.Range("N14").Value = 0.9245046239218979
$ws.Range("O14").Value = 1.916695798392539
$ws.Range("B15").Value = 1.053565251177986
$ws.Range("C15").Value = 0.237535242794678
$ws.Range("D15").Value = 0.01823185873121957
$ws.Range("E15").Value = 0.1592142078853627
$ws.Range("F15").Value = 0.6081859078227865
$ws.Range("M15").Value = 0.4699901019325665
$ws.Range("N15").Value = 0.9256953409738102
$ws.Range("O15").Value = 1.91285207802045
$ws.Range("B16").Value = 0.9903354920342053
$ws.Range("C16").Value = 0.2244841974924441
$ws.Range("D16").Value = 0.01736485189015724
$ws.Range("E16").Value = 0.1536043294146197
$ws.Range("F16").Value = 0.5984610503702896
$ws.Range("M16").Value = 0.4433958892782783
$ws.Range("N16").Value = 0.9326549809655376
$ws.Range("O16").Value = 1.891277687039491
$ws.Range("B17").Value = 0.9515313190295274
$ws.Range("C17").Value = 0.2164604706176476
$ws.Range("D17").Value = 0.01683167659071216
$ws.Range("E17").Value = 0.1501950311698508
$ws.Range("F17").Value = 0.5926170893340412
$ws.Range("M17").Value = 0.4271021592477524
$ws.Range("N17").Value = 0.9370458749670192
$ws.Range("O17").Value = 1.878445309108656
$ws.Range("B18").Value = 0.9292070599042859
$ws.Range("C18").Value = 0.211839132781904
$ws.Range("D18").Value = 0.01652453964079825
$ws.Range("E18").Value = 0.1482458291065782
$ws.Range("F18").Value = 0.5893005104469751
$ws.Range("M18").Value = 0.4177381277398666
$ws.Range("N18").Value = 0.9396159752171585
$ws.Range("O18").Value = 1.871212144837898
$ws.Range("B19").Value = 0.9216475976369907
$ws.Range("C19").Value = 0.2102733511798931
$ws.Range("D19").Value = 0.01642046854230017
$ws.Range("E19").Value = 0.1475878699328277
$ws.Range("F19").Value = 0.5881852444154561
$ws.Range("M19").Value = 0.4145689474829055
$ws.Range("N19").Value = 0.9404938196128754
$ws.Range("O19").Value = 1.868788450514728
$ws.Range("B20").Value = 0.9556626267054185
$ws.Range("C20").Value = 0.2173152638379747
$ws.Range("D20").Value = 0.01688848267323095
$ws.Range("E20").Value = 0.1505567401156895
$ws.Range("F20").Value = 0.5932345592257775
$ws.Range("M20").Value = 0.4288358583875578
$ws.Range("N20").Value = 0.936573842815541
$ws.Range("O20").Value = 1.879796045364458
$ws.Range("B21").Value = 1.069885123073732
$ws.Range("C21").Value = 0.2408993523890217
$ws.Range("D21").Value = 0.01845529675152591
$ws.Range("E21").Value = 0.1606725907553965
$ws.Range("F21").Value = 0.6107344452886707
$ws.Range("M21").Value = 0.476862693354775
$ws.Range("N21").Value = 0.9239360518292301
$ws.Range("O21").Value = 1.918546982150332
$ws.Range("B22").Value = 1.14448411212328
$ws.Range("C22").Value = 0.2562553624709665
$ws.Range("D22").Value = 0.01947497517326013
$ws.Range("E22").Value = 0.1673903097265494
$ws.Range("F22").Value = 0.6225717856963229
$ws.Range("M22").Value = 0.5083197299137368
$ws.Range("N22").Value = 0.9160741583198018
$ws.Range("O22").Value = 1.945195510101996
$ws.Range("B23").Value = 1.104674295187863
$ws.Range("C23").Value = 0.2480649063973601
$ws.Range("D23").Value = 0.01893115670711154
$ws.Range("E23").Value = 0.1637951056823113
$ws.Range("F23").Value = 0.6162172640932226
$ws.Range("M23").Value = 0.4915242176267611
$ws.Range("N23").Value = 0.9202337802595864
$ws.Range("O23").Value = 1.930851313903275
$ws.Range("B24").Value = 0.9537949098552758
$ws.Range("C24").Value = 0.2169288379442946
$ws.Range("D24").Value = 0.01686280253274219
$ws.Range("E24").Value = 0.1503931777315017
$ws.Range("F24").Value = 0.5929552665679267
$ws.Range("M24").Value = 0.4280520422652643
$ws.Range("N24").Value = 0.9367871061485076
$ws.Range("O24").Value = 1.87918492796652
$ws.Range("B25").Value = 0.7908543311913832
$ws.Range("C25").Value = 0.1831033339350938
$ws.Range("D25").Value = 0.01461391480474816
$ws.Range("E25").Value = 0.1363845327678703
$ws.Range("F25").Value = 0.5695715863173945
$ws.Range("M25").Value = 0.359880677049226
$ws.Range("N25").Value = 0.9563463693887044
$ws.Range("O25").Value = 1.829095709864788
